$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for the Price/Volume columns (and the B/C columns
# touched on row 51) so numeric-looking strings are not coerced into numbers.
$textRange = $ws.Range("B2:E51")
$textRange.NumberFormat = "@"

$ws.Range('D2').Value = '65.346.17'
$ws.Range('E2').Value = '  -0.63%  '
$ws.Range('D3').Value = '2.952.38'
$ws.Range('E3').Value = '  -1.67%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '569.41'
$ws.Range('E5').Value = '  -2.72%  '
$ws.Range('D6').Value = '159.97'
$ws.Range('E6').Value = '  +4.02%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +0.60%  '
$ws.Range('D9').Value = '2.945.39'
$ws.Range('E9').Value = '  -1.77%  '
$ws.Range('D10').Value = '6.65'
$ws.Range('E10').Value = '  -4.69%  '
$ws.Range('E11').Value = '  -1.46%  '
$ws.Range('D12').Value = '0.458'
$ws.Range('E12').Value = '  +2.37%  '
$ws.Range('E13').Value = '  +2.08%  '
$ws.Range('D14').Value = '34.32'
$ws.Range('E14').Value = '  +0.84%  '
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').Value = '65.306.82'
$ws.Range('E16').Value = '  -0.52%  '
$ws.Range('D17').Value = '3.439.72'
$ws.Range('E17').Value = '  -1.67%  '
$ws.Range('D18').Value = '6.97'
$ws.Range('E18').Value = '  +0.65%  '
$ws.Range('D19').Value = '2.951.08'
$ws.Range('E19').Value = '  -1.67%  '
$ws.Range('D20').Value = '14.61'
$ws.Range('E20').Value = '  +6.58%  '
$ws.Range('E21').Value = '  -1.78%  '
$ws.Range('D22').Value = '0.689'
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('D23').Value = '7.26'
$ws.Range('E23').Value = '  -0.99%  '
$ws.Range('D24').Value = '82.29'
$ws.Range('E24').Value = '  +1.09%  '
$ws.Range('E25').Value = '  -0.74%  '
$ws.Range('D26').Value = '12.12'
$ws.Range('D27').Value = '10.04'
$ws.Range('E27').Value = '  -5.80%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').Value = '8.06'
$ws.Range('E29').Value = '  +3.72%  '
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('E32').Value = '  -1.95%  '
$ws.Range('D33').Value = '27.14'
$ws.Range('E33').Value = '  +0.78%  '
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('E36').Value = '  -1.52%  '
$ws.Range('D37').Value = '5.72'
$ws.Range('E37').Value = '  -1.01%  '
$ws.Range('D38').Value = '49.01'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').Value = '44.23'
$ws.Range('E39').Value = '  -3.12%  '
$ws.Range('D40').Value = '1.98'
$ws.Range('E40').Value = '  -6.68%  '
$ws.Range('D41').Value = '2.86'
$ws.Range('E41').Value = '  -1.86%  '
$ws.Range('E42').Value = '  -1.30%  '
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('D44').Value = '8.48'
$ws.Range('E44').Value = '  +0.75%  '
$ws.Range('D45').Value = '384.44'
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('E46').Value = '  +0.23%  '
$ws.Range('D47').Value = '2.716.37'
$ws.Range('E47').Value = '  -1.82%  '
$ws.Range('D48').Value = '133.20'
$ws.Range('E48').Value = '  -1.55%  '
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('E50').Value = '  +5.53%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '23.31'
$ws.Range('E51').Value = '  -0.07%  '

# Restore the default "Normal" style so we do not leave a stray text format behind.
$textRange.Style = "Normal"
